$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (d=2)
$ws.Range("G2").Value = 849585974.7623143
$ws.Range("H2").Value = 46128147165.71857
$ws.Range("I2").Value = 1177234628194075

# Row 3 (d=3)
$ws.Range("B3").Value = 0.7074335935529131
$ws.Range("C3").Value = 29.2968893301947
$ws.Range("D3").Value = 515.8738629160657
$ws.Range("E3").Value = 50728.35286098456
$ws.Range("F3").Value = 4652938.886985654
$ws.Range("G3").Value = 241300287.4162372
$ws.Range("H3").Value = 14711932116.51295
$ws.Range("I3").Value = 486637300837981.1

# Row 4 (d=4)
$ws.Range("B4").Value = 0.6994382448879282
$ws.Range("C4").Value = 34.413495310372
$ws.Range("D4").Value = 1822.124109082956
$ws.Range("E4").Value = 211628.9445684948
$ws.Range("F4").Value = 19664540.92785762
$ws.Range("G4").Value = 878281275.2445639
$ws.Range("H4").Value = 47381874240.55646
$ws.Range("I4").Value = 1195642070166485

# Row 5 (d=5)
$ws.Range("G5").Value = 796011138.9636092
$ws.Range("H5").Value = 43858732466.16466
$ws.Range("I5").Value = 1153203812976251

# Row 6 (d=6)
$ws.Range("B6").Value = 0.7274922139810557
$ws.Range("C6").Value = 30.54165270518868
$ws.Range("D6").Value = 804.4784241064688
$ws.Range("E6").Value = 67650.91406843558
$ws.Range("F6").Value = 5925536.849115291
$ws.Range("G6").Value = 290279923.7108551
$ws.Range("H6").Value = 16350848539.72966
$ws.Range("I6").Value = 423202772033612.8

# Row 7 (d=7)
$ws.Range("G7").Value = 25243545.10549158
$ws.Range("H7").Value = 1541935848.40369
$ws.Range("I7").Value = 54869063329978.05

# Row 8 (d=10)
$ws.Range("G8").Value = 160826311.8691058
$ws.Range("H8").Value = 9982603070.527473
$ws.Range("I8").Value = 340060298262874.4

# Row 9 (AREPD)
$ws.Range("B9").Value = 0.6889183596399426
$ws.Range("C9").Value = 29.35904219599299
$ws.Range("D9").Value = 597.7100550895532
$ws.Range("E9").Value = 85313.40533540784
$ws.Range("F9").Value = 9165618.36345561
$ws.Range("G9").Value = 455451865.039818
$ws.Range("H9").Value = 27058168820.91214
$ws.Range("I9").Value = 817390885859121.2

# Row 10 (AV-MCPS)
$ws.Range("G10").Value = 26.63349199135783
$ws.Range("H10").Value = 1570.67474124423
$ws.Range("I10").Value = 67958618.98401152
